$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 223 (existing data rows 223-298 shift
# down to become rows 225-300), mirroring a new week ("Primera" + "Segunda")
# of Mango price entries being logged at the top of the data block.
$ws.Rows.Item(223).Insert()
$ws.Rows.Item(223).Insert()

# New row 223 - Mango, "Primera" quality, week of 2021-12-29
$ws.Cells.Item(223, 1).Value = 3
$ws.Cells.Item(223, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(223, 3).Value = "Coquimbo"
$ws.Cells.Item(223, 4).Value = 44559
$ws.Cells.Item(223, 5).Value = 5
$ws.Cells.Item(223, 6).Value = "Fruta"
$ws.Cells.Item(223, 7).Value = 100108
$ws.Cells.Item(223, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(223, 9).Value = 100108002
$ws.Cells.Item(223, 10).Value = "Mango"
$ws.Cells.Item(223, 11).Value = "Sin especificar"
$ws.Cells.Item(223, 12).Value = "Primera"
$ws.Cells.Item(223, 13).Value = 456
$ws.Cells.Item(223, 14).Value = 7000
$ws.Cells.Item(223, 15).Value = 7000
$ws.Cells.Item(223, 16).Value = 7000
$ws.Cells.Item(223, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(223, 18).Value = "Perú"
$ws.Cells.Item(223, 19).Value = 1750
$ws.Cells.Item(223, 20).Value = 4

# New row 224 - Mango, "Segunda" quality, same week (2021-12-29)
$ws.Cells.Item(224, 1).Value = 3
$ws.Cells.Item(224, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(224, 3).Value = "Coquimbo"
$ws.Cells.Item(224, 4).Value = 44559
$ws.Cells.Item(224, 5).Value = 5
$ws.Cells.Item(224, 6).Value = "Fruta"
$ws.Cells.Item(224, 7).Value = 100108
$ws.Cells.Item(224, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(224, 9).Value = 100108002
$ws.Cells.Item(224, 10).Value = "Mango"
$ws.Cells.Item(224, 11).Value = "Sin especificar"
$ws.Cells.Item(224, 12).Value = "Segunda"
$ws.Cells.Item(224, 13).Value = 456
$ws.Cells.Item(224, 14).Value = 7000
$ws.Cells.Item(224, 15).Value = 7000
$ws.Cells.Item(224, 16).Value = 7000
$ws.Cells.Item(224, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(224, 18).Value = "Perú"
$ws.Cells.Item(224, 19).Value = 1750
$ws.Cells.Item(224, 20).Value = 4

Write-Output "rows inserted and populated"
